$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118, shifting existing rows 118-142 down to 119-143
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new data point
$ws.Cells.Item(118, 1).Value = 7
$ws.Cells.Item(118, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(118, 3).Value = "Ñuble"
$ws.Cells.Item(118, 4).Value = 44505
$ws.Cells.Item(118, 4).NumberFormat = $ws.Cells.Item(119, 4).NumberFormat
$ws.Cells.Item(118, 5).Value = 16
$ws.Cells.Item(118, 6).Value = 100112006
$ws.Cells.Item(118, 7).Value = "Repollo"
$ws.Cells.Item(118, 8).Value = "Crespo record"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 400
$ws.Cells.Item(118, 11).Value = 600
$ws.Cells.Item(118, 12).Value = 700
$ws.Cells.Item(118, 13).Value = 650
$ws.Cells.Item(118, 14).Value = "$/unidad"
$ws.Cells.Item(118, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(118, 16).Value = 650
$ws.Cells.Item(118, 17).Value = 1
$ws.Cells.Item(118, 18).Value = "Hortaliza"
